$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, centered, bordered) onto the new I1/J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF), rows 2-27
$data = @(
    @(2, 1, 6),
    @(3, 1, 5),
    @(4, 2, 5),
    @(5, 1, 5),
    @(6, 1, 5),
    @(7, 2, 4),
    @(8, 1, 4),
    @(9, 1, 4),
    @(10, 1, 6),
    @(11, 1, 5),
    @(12, 1, 6),
    @(13, 2, 7),
    @(14, 1, 5),
    @(15, 1, 5),
    @(16, 2, 7),
    @(17, 1, 5),
    @(18, 2, 6),
    @(19, 1, 3),
    @(20, 2, 7),
    @(21, 1, 6),
    @(22, 1, 6),
    @(23, 1, 6),
    @(24, 1, 5),
    @(25, 8, 9),
    @(26, 3, 4),
    @(27, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
